$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8957854509353638
$ws.Range("B1").Value = 1.785931944847107
$ws.Range("C1").Value = 4.2351975440979
$ws.Range("D1").Value = 3.504871368408203
$ws.Range("E1").Value = 1.504392743110657
